$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_7_4_0"
$ws.Range("B2").Value = 0.7010364209363523
$ws.Range("C2").Value = 0.9953471598035962
$ws.Range("D2").Value = -0.02903659871527808
$ws.Range("E2").Value = 0.9670915577989236
$ws.Range("F2").Value = 0.3308645784854889
$ws.Range("G2").Value = 0.01174833066761494
$ws.Range("H2").Value = 0.09948210418224335
$ws.Range("I2").Value = 0.05303477123379707
$ws.Range("A3").Value = "model_7_4_1"
$ws.Range("B3").Value = 0.7116203156363476
$ws.Range("C3").Value = 0.994811192539603
$ws.Range("D3").Value = -0.06336989484612077
$ws.Range("E3").Value = 0.9656776107650619
$ws.Range("F3").Value = 0.3191513419151306
$ws.Range("G3").Value = 0.01310163643211126
$ws.Range("H3").Value = 0.1028012707829475
$ws.Range("I3").Value = 0.05531346425414085
$ws.Range("A4").Value = "model_7_4_2"
$ws.Range("B4").Value = 0.7215422290425451
$ws.Range("C4").Value = 0.9942220113491079
$ws.Range("D4").Value = -0.1167933056692416
$ws.Range("E4").Value = 0.9636809891137257
$ws.Range("F4").Value = 0.3081707060337067
$ws.Range("G4").Value = 0.01458930782973766
$ws.Range("H4").Value = 0.1079659759998322
$ws.Range("I4").Value = 0.05853119120001793
$ws.Range("A5").Value = "model_7_4_3"
$ws.Range("B5").Value = 0.7314160228038955
$ws.Range("C5").Value = 0.9937581363616567
$ws.Range("D5").Value = -0.1433477429565704
$ws.Range("E5").Value = 0.9625465148118831
$ws.Range("F5").Value = 0.2972433269023895
$ws.Range("G5").Value = 0.01576058380305767
$ws.Range("H5").Value = 0.110533133149147
$ws.Range("I5").Value = 0.06035949289798737
$ws.Range("A6").Value = "model_7_4_4"
$ws.Range("B6").Value = 0.7410361528003662
$ws.Range("C6").Value = 0.9933381818486436
$ws.Range("D6").Value = -0.1579419036743737
$ws.Range("E6").Value = 0.9617861650210743
$ws.Range("F6").Value = 0.2865966856479645
$ws.Range("G6").Value = 0.01682095974683762
$ws.Range("H6").Value = 0.1119440197944641
$ws.Range("I6").Value = 0.06158486381173134
$ws.Range("A7").Value = "model_7_4_5"
$ws.Range("B7").Value = 0.7503608011788359
$ws.Range("C7").Value = 0.9929281549487654
$ws.Range("D7").Value = -0.1654643238924884
$ws.Range("E7").Value = 0.961233911593399
$ws.Range("F7").Value = 0.2762770354747772
$ws.Range("G7").Value = 0.01785627007484436
$ws.Range("H7").Value = 0.1126712560653687
$ws.Range("I7").Value = 0.06247486919164658
$ws.Range("A8").Value = "model_7_4_6"
$ws.Range("B8").Value = 0.7593124892612211
$ws.Range("C8").Value = 0.9924396026849251
$ws.Range("D8").Value = -0.1706728013615191
$ws.Range("E8").Value = 0.9606813526504161
$ws.Range("F8").Value = 0.2663701176643372
$ws.Range("G8").Value = 0.01908985525369644
$ws.Range("H8").Value = 0.11317478120327
$ws.Range("I8").Value = 0.06336536258459091
$ws.Range("A9").Value = "model_7_4_7"
$ws.Range("B9").Value = 0.7679069144789443
$ws.Range("C9").Value = 0.9918133779513236
$ws.Range("D9").Value = -0.1721146476716549
$ws.Range("E9").Value = 0.9601214167747684
$ws.Range("F9").Value = 0.2568586468696594
$ws.Range("G9").Value = 0.02067106030881405
$ws.Range("H9").Value = 0.1133141741156578
$ws.Range("I9").Value = 0.06426775455474854
$ws.Range("A10").Value = "model_7_4_8"
$ws.Range("B10").Value = 0.7761133247617075
$ws.Range("C10").Value = 0.9910541349790934
$ws.Range("D10").Value = -0.1722857151106474
$ws.Range("E10").Value = 0.9594868893542909
$ws.Range("F10").Value = 0.2477765530347824
$ws.Range("G10").Value = 0.02258813381195068
$ws.Range("H10").Value = 0.113330714404583
$ws.Range("I10").Value = 0.06529034674167633
$ws.Range("A11").Value = "model_7_4_9"
$ws.Range("B11").Value = 0.7839656052351758
$ws.Range("C11").Value = 0.9901941677396083
$ws.Range("D11").Value = -0.1708206771694321
$ws.Range("E11").Value = 0.9588147970333097
$ws.Range("F11").Value = 0.2390864044427872
$ws.Range("G11").Value = 0.02475953474640846
$ws.Range("H11").Value = 0.1131890788674355
$ws.Range("I11").Value = 0.06637348234653473
$ws.Range("A12").Value = "model_7_4_10"
$ws.Range("B12").Value = 0.7915069777818194
$ws.Range("C12").Value = 0.9892221786658402
$ws.Range("D12").Value = -0.1654686305889685
$ws.Range("E12").Value = 0.9581598176206432
$ws.Range("F12").Value = 0.2307403385639191
$ws.Range("G12").Value = 0.02721378579735756
$ws.Range("H12").Value = 0.1126716732978821
$ws.Range("I12").Value = 0.06742902845144272
$ws.Range("A13").Value = "model_7_4_11"
$ws.Range("B13").Value = 0.7987333832400106
$ws.Range("C13").Value = 0.9880667076645869
$ws.Range("D13").Value = -0.1557372249063298
$ws.Range("E13").Value = 0.9574758523560566
$ws.Range("F13").Value = 0.2227428406476974
$ws.Range("G13").Value = 0.03013132698833942
$ws.Range("H13").Value = 0.1117308884859085
$ws.Range("I13").Value = 0.06853130459785461
$ws.Range("A14").Value = "model_7_4_12"
$ws.Range("B14").Value = 0.8055656477984454
$ws.Range("C14").Value = 0.9866314897918743
$ws.Range("D14").Value = -0.1436969597742572
$ws.Range("E14").Value = 0.9566251897023411
$ws.Range("F14").Value = 0.2151815146207809
$ws.Range("G14").Value = 0.03375522419810295
$ws.Range("H14").Value = 0.1105668917298317
$ws.Range("I14").Value = 0.0699022188782692
$ws.Range("A15").Value = "model_7_4_13"
$ws.Range("B15").Value = 0.8120296626199467
$ws.Range("C15").Value = 0.9848955053489047
$ws.Range("D15").Value = -0.1298804371649824
$ws.Range("E15").Value = 0.9555755315290263
$ws.Range("F15").Value = 0.2080277651548386
$ws.Range("G15").Value = 0.03813855350017548
$ws.Range("H15").Value = 0.1092311814427376
$ws.Range("I15").Value = 0.07159382849931717
$ws.Range("A16").Value = "model_7_4_14"
$ws.Range("B16").Value = 0.8180600932885804
$ws.Range("C16").Value = 0.9827285379919948
$ws.Range("D16").Value = -0.1155430589260462
$ws.Range("E16").Value = 0.954182708291053
$ws.Range("F16").Value = 0.2013538628816605
$ws.Range("G16").Value = 0.0436101034283638
$ws.Range("H16").Value = 0.1078451126813889
$ws.Range("I16").Value = 0.0738384798169136
$ws.Range("A17").Value = "model_7_4_15"
$ws.Range("B17").Value = 0.8237079162181447
$ws.Range("C17").Value = 0.9802214206070941
$ws.Range("D17").Value = -0.09969650401899832
$ws.Range("E17").Value = 0.9525505685244295
$ws.Range("F17").Value = 0.1951033920049667
$ws.Range("G17").Value = 0.04994052276015282
$ws.Range("H17").Value = 0.1063131541013718
$ws.Range("I17").Value = 0.0764688178896904
$ws.Range("A18").Value = "model_7_4_16"
$ws.Range("B18").Value = 0.8290097938278064
$ws.Range("C18").Value = 0.9773600867197085
$ws.Range("D18").Value = -0.08184012228103632
$ws.Range("E18").Value = 0.9506814332620406
$ws.Range("F18").Value = 0.1892357766628265
$ws.Range("G18").Value = 0.05716533958911896
$ws.Range("H18").Value = 0.1045868843793869
$ws.Range("I18").Value = 0.07948108017444611
$ws.Range("A19").Value = "model_7_4_17"
$ws.Range("B19").Value = 0.833956291664372
$ws.Range("C19").Value = 0.9742006670297046
$ws.Range("D19").Value = -0.06278338586265719
$ws.Range("E19").Value = 0.94859856023617
$ws.Range("F19").Value = 0.1837614625692368
$ws.Range("G19").Value = 0.06514281034469604
$ws.Range("H19").Value = 0.1027445793151855
$ws.Range("I19").Value = 0.08283782005310059
$ws.Range("A20").Value = "model_7_4_18"
$ws.Range("B20").Value = 0.8386455507136568
$ws.Range("C20").Value = 0.9709915185074701
$ws.Range("D20").Value = -0.04415397477127847
$ws.Range("E20").Value = 0.946462796452093
$ws.Range("F20").Value = 0.1785718202590942
$ws.Range("G20").Value = 0.07324584573507309
$ws.Range("H20").Value = 0.1009435728192329
$ws.Range("I20").Value = 0.08627977967262268
$ws.Range("A21").Value = "model_7_4_19"
$ws.Range("B21").Value = 0.8431312028730571
$ws.Range("C21").Value = 0.9677719315401913
$ws.Range("D21").Value = -0.02494644461521722
$ws.Range("E21").Value = 0.944334302443721
$ws.Range("F21").Value = 0.1736075133085251
$ws.Range("G21").Value = 0.08137524127960205
$ws.Range("H21").Value = 0.09908668696880341
$ws.Range("I21").Value = 0.0897100418806076
$ws.Range("A22").Value = "model_7_4_20"
$ws.Range("B22").Value = 0.8473729622258691
$ws.Range("C22").Value = 0.9644399825243707
$ws.Range("D22").Value = -0.005297002176353649
$ws.Range("E22").Value = 0.9421255319024705
$ws.Range("F22").Value = 0.1689131557941437
$ws.Range("G22").Value = 0.08978834748268127
$ws.Range("H22").Value = 0.09718708693981171
$ws.Range("I22").Value = 0.09326965361833572
$ws.Range("A23").Value = "model_7_4_21"
$ws.Range("B23").Value = 0.851405486088011
$ws.Range("C23").Value = 0.961088198580199
$ws.Range("D23").Value = 0.01501104771504935
$ws.Range("E23").Value = 0.9399184975149356
$ws.Range("F23").Value = 0.1644503325223923
$ws.Range("G23").Value = 0.09825152903795242
$ws.Range("H23").Value = 0.0952237993478775
$ws.Range("I23").Value = 0.09682648628950119
$ws.Range("A24").Value = "model_7_4_22"
$ws.Range("B24").Value = 0.8552787805371339
$ws.Range("C24").Value = 0.9578441704621208
$ws.Range("D24").Value = 0.03385332141844122
$ws.Range("E24").Value = 0.9377594922033478
$ws.Range("F24").Value = 0.1601637452840805
$ws.Range("G24").Value = 0.1064426451921463
$ws.Range("H24").Value = 0.09340222179889679
$ws.Range("I24").Value = 0.1003059074282646
$ws.Range("A25").Value = "model_7_4_23"
$ws.Range("B25").Value = 0.858953942250995
$ws.Range("C25").Value = 0.954584552295627
$ws.Range("D25").Value = 0.05273664968149294
$ws.Range("E25").Value = 0.9355886252203311
$ws.Range("F25").Value = 0.1560964286327362
$ws.Range("G25").Value = 0.1146731078624725
$ws.Range("H25").Value = 0.09157668054103851
$ws.Range("I25").Value = 0.1038044318556786
$ws.Range("A26").Value = "model_7_4_24"
$ws.Range("B26").Value = 0.8624568625104797
$ws.Range("C26").Value = 0.9513572556322939
$ws.Range("D26").Value = 0.07117923348635713
$ws.Range("E26").Value = 0.933432494172218
$ws.Range("F26").Value = 0.1522197425365448
$ws.Range("G26").Value = 0.1228219792246819
$ws.Range("H26").Value = 0.08979374170303345
$ws.Range("I26").Value = 0.1072792261838913